$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Title paragraph: drop the stray "_GoBack" bookmark that was
#    sitting around the title run.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Body paragraph: the author's cursor (and hence the new "_GoBack"
#    bookmark) ended up mid-word, between "...924 is the ratin" and
#    "g for season one...". Recreate that split.
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("924 is the ratin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(0)
    $r.Bookmarks.Add("_GoBack")
}

# ------------------------------------------------------------------
# 3) Footer: "A-Diakoptics suite for OpenDSS" becomes
#    "Seasonal rating simulation in OpenDSS" (matching the new title),
#    with "OpenDSS" kept as its own run.
# ------------------------------------------------------------------
$sec = $d.Sections(1)
$footer = $sec.Footers(1)

$fr = $footer.Range
$fr.Find.Execute("A-Diakoptics suite for OpenDSS", $true, $false, $false, $false, $false, $true, 1, $false, "Seasonal rating simulation in OpenDSS", 2) | Out-Null

$fr2 = $footer.Range
$found2 = $fr2.Find.Execute("OpenDSS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    # Toggling a character property on just this word forces Word to
    # split it into its own run, then we put the property back so the
    # formatting is unchanged from the rest of the line.
    $fr2.Font.Bold = $true
    $fr2.Font.Bold = $false
}
